# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2
# of the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 11:22:29"
$wsZhCn.Range("H2").Value = "2016-03-17 11:22:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 11:22:32"
$wsDeDe.Range("H2").Value = "2016-03-17 11:22:55"
